# Thai translation pass for "Email 10-1 [TEMPLATE] Partner email - thank you
# email (without photos).docx".
#
# Word constants used below (spelled out as literals since this is a
# constant-free COM-interop snippet):
#   wdReplaceOne = 1
#   wdReplaceAll = 2
#   wdFindContinue = 1 (Wrap parameter)

$d = $word.ActiveDocument

# --- Paragraph 1: "English / Portuguese / French / Thai / Vietnamese / Spanish" header ---
# The word "English" lives inside a hyperlink run; replace it on its own so we
# don't touch the following run's formatting.
$p = $d.Paragraphs.Item(1).Range
$p.Find.Execute("English", $true, $false, $false, $false, $false, $true, 1, $false, "ภาษาอังกฤษ", 1) | Out-Null

# The rest of the line (" / Portuguese / French / Thai / Vietnamese / Spanish")
# is a separate run right after the hyperlink run. Search starting from
# "Portuguese" (i.e. skip the leading " / ", which is unchanged) so the match
# does not begin exactly on the hyperlink/plain-text run boundary.
$p = $d.Paragraphs.Item(1).Range
$p.Find.Execute("Portuguese / French / Thai / Vietnamese / Spanish", $true, $false, $false, $false, $false, $true, 1, $false, "ภาษาโปรตุเกส / ภาษาฝรั่งเศส /ภาษาไทย / ภาษาเวียดนาม / ภาษาสเปน", 1) | Out-Null

# --- Paragraph 3: lone "English" label ---
$p = $d.Paragraphs.Item(3).Range
$p.Find.Execute("English", $true, $false, $false, $false, $false, $true, 1, $false, "ภาษาอังกฤษ", 1) | Out-Null

# --- Paragraph 5: "Brief" label ---
$p = $d.Paragraphs.Item(5).Range
$p.Find.Execute("Brief", $true, $false, $false, $false, $false, $true, 1, $false, "บทย่อ", 1) | Out-Null

# --- Paragraph 8: "Target audience" label ---
$p = $d.Paragraphs.Item(8).Range
$p.Find.Execute("Target audience", $true, $false, $false, $false, $false, $true, 1, $false, "กลุ่มเป้าหมาย", 1) | Out-Null

# --- Paragraph 13: "You made our event a success! <emoji>" heading ---
# Leave the trailing party-popper emoji glyph untouched.
$p = $d.Paragraphs.Item(13).Range
$p.Find.Execute("You made our event a success!", $true, $false, $false, $false, $false, $true, 1, $false, "คุณได้ช่วยให้กิจกรรมของเราประสบความสำเร็จด้วยดี!", 1) | Out-Null

# --- Paragraph 15: "Hi [PARTNER NAME], " ---
$p = $d.Paragraphs.Item(15).Range
$p.Find.Execute("Hi ", $true, $false, $false, $false, $false, $true, 1, $false, "สวัสดี ", 1) | Out-Null

$p = $d.Paragraphs.Item(15).Range
$p.Find.Execute(", ", $true, $false, $false, $false, $false, $true, 1, $false, " ", 1) | Out-Null

# --- Paragraph 17: "Thank you for attending [EVENT NAME] in [CITY], [COUNTRY]. We hope ..." ---
$p = $d.Paragraphs.Item(17).Range
$p.Find.Execute("Thank you for attending ", $true, $false, $false, $false, $false, $true, 1, $false, "ขอบคุณที่ได้เข้าร่วมงาน ", 1) | Out-Null

$p = $d.Paragraphs.Item(17).Range
$p.Find.Execute(" in ", $true, $false, $false, $false, $false, $true, 1, $false, " ซึ่งจัดขึ้นที่เมือง ", 1) | Out-Null

$p = $d.Paragraphs.Item(17).Range
$p.Find.Execute(", ", $true, $false, $false, $false, $false, $true, 1, $false, " ประเทศ ", 1) | Out-Null

$p = $d.Paragraphs.Item(17).Range
$p.Find.Execute(". We hope you had a great time, and it was a pleasure getting to know you!", $true, $false, $false, $false, $false, $true, 1, $false, " เราหวังว่า คุณได้มีช่วงเวลาที่ดีเยี่ยม และเราก็มีความยินดีอย่างยิ่งที่ได้รู้จักกับคุณ!", 1) | Out-Null

# --- Paragraph 27: closing line ---
$p = $d.Paragraphs.Item(27).Range
$p.Find.Execute("We hope the event inspired you as much as it did us, and let" + [char]0x2019 + "s keep growing together!", $true, $false, $false, $false, $false, $true, 1, $false, "เราหวังว่า กิจกรรมนี้ได้สร้างแรงบันดาลใจให้คุณมากพอๆ กับที่ได้สร้างให้กับเรา แล้วเรามาเติบโตไปด้วยกันเถอะ!", 1) | Out-Null

# --- Comment text: "choose either one" ---
$c = $d.Comments.Item(1)
$c.Range.Text = "เลือกอย่างใดอย่างหนึ่ง"
